$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Final layout (per the diff):
#   Row1 headers: Nombre | Apellido | Run        | Correo
#   Row3:         Luis Miguel | Leon Valenzuela | 17428518-7 | lmleon@outlook.com
#   Row4:         Rojas       | Rojas Arias     | 17417030-4 | b.rojas@icci.cl
# ---------------------------------------------------------------------------

# ---- 1. Update existing row 3, and create row 4 with the new member ----
$ws.Range("A3").Value = "Luis Miguel"
$ws.Range("A4").Value = "Rojas"
$ws.Range("B4").Value = "Rojas Arias"
$ws.Range("C4").Value = "17417030-4"

# ---- 2. Expand the table FIRST (creates placeholder 4th column) ----
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:D7"))

# ---- 3. Now set the header text - this renames the table column too ----
$ws.Range("D1").Value = "Correo"

# ---- 4. Add the e-mail values and turn them into mailto hyperlinks ----
$ws.Range("D3").Value = "lmleon@outlook.com"
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:lmleon@outlook.com")
$ws.Range("D4").Value = "b.rojas@icci.cl"
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:b.rojas@icci.cl")

# ---- 5. Leftover hyperlink-styled empty cell (matches author's original edit) ----
$ws.Range("D5").Style = "Hyperlink"

# ---- 6. Column width for the new column ----
$ws.Columns.Item(4).AutoFit()

# ---- 7. Selection, matching the author's final cursor position ----
$ws.Range("C9").Select()
